$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'25.769.54"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -3.97%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.818.50"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -2.90%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.23%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'279.62"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -7.25%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'1.0000"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Value = "'0.5115"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -4.27%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.3545"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -5.35%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'44.61"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -1.99%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.06669"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -7.34%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'20.11"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -6.97%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.8284"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -6.94%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.07906"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -3.40%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'1.820.93"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -2.87%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'5.087"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -4.15%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'88.03"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -5.56%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'1.0000"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -0.30%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'14.10"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -5.00%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.000008041"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -5.78%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'  -0.16%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'25.808.44"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Value = "'4.759"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -4.64%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'9.979"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -6.02%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'6.128"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -3.95%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'2.241"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -1.80%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'142.39"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -2.57%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  -4.08%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'17.18"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -4.94%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'109.40"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -4.03%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  -7.93%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'4.244"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -8.06%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'0.08771"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -3.70%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'0.04913"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -2.01%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'0.7329"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -9.63%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'1.141"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -2.75%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  -2.95%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'3.159"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -1.46%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'2.391"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -9.76%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.01855"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -5.23%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.5169"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -15.12%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.9666"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -9.60%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'6.248"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -5.12%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'111.26"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -3.30%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'8.067"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -9.03%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'  -0.20%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.4582"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -11.24%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.1372"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -8.27%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'36.66"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -2.28%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'9.221"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -7.27%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'1.505"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -8.32%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.05816"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -3.98%  "
$ws.Range("E51").Style = "Normal"
